$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "season record" header columns (Wins/Losses/Ties),
# reusing the existing bordered/bold header style from column AC.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate each player's row with the team's season record.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 98
    $ws.Cells.Item($r, 31).Value = 64
    $ws.Cells.Item($r, 32).Value = 0
}
